$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated coin values (prices/volume%) scraped on 3-1-2023 run.
# Columns B/C (coin name / link) are plain text; columns D/E (price / volume%)
# are numeric-looking strings stored as TEXT in the source file, so we prefix
# them with a leading apostrophe to force Excel to keep them as text, matching
# the original inlineStr cell type instead of letting Excel coerce to a number.

$ws.Range("D2").Value = "'245.33"
$ws.Range("E2").Value = "'-0.62%"
$ws.Range("D3").Value = "'28.41"
$ws.Range("E3").Value = "'-4.39%"
$ws.Range("D4").Value = "'5.248"
$ws.Range("E4").Value = "'1.60%"
$ws.Range("E5").Value = "'-0.47%"
$ws.Range("D6").Value = "'6.633"
$ws.Range("E6").Value = "'0.67%"
$ws.Range("D7").Value = "'3.202"
$ws.Range("E7").Value = "'3.37%"
$ws.Range("D8").Value = "'0.8507"
$ws.Range("E8").Value = "'-0.69%"
$ws.Range("D9").Value = "'0.9050"
$ws.Range("E9").Value = "'4.39%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1368"
$ws.Range("E10").Value = "'0.20%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07077"
$ws.Range("E11").Value = "'0.17%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03154"
$ws.Range("E12").Value = "'7.60%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09192"
$ws.Range("E13").Value = "'-2.01%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001531"
$ws.Range("E14").Value = "'0.99%"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").Value = "'0.0005986"
$ws.Range("E15").Value = "'-0.17%"
$ws.Range("D16").Value = "'0.005955"
$ws.Range("E16").Value = "'-2.54%"
$ws.Range("D17").Value = "'3.489"
$ws.Range("E17").Value = "'-0.01%"
$ws.Range("E18").Value = "'-0.60%"
$ws.Range("D19").Value = "'0.3171"
$ws.Range("E19").Value = "'-0.37%"
$ws.Range("D20").Value = "'0.03306"
$ws.Range("E20").Value = "'-2.55%"
$ws.Range("D21").Value = "'0.1283"
$ws.Range("E21").Value = "'-1.34%"
$ws.Range("D22").Value = "'3.524"
$ws.Range("E22").Value = "'1.46%"
$ws.Range("D23").Value = "'0.04076"
$ws.Range("E23").Value = "'-1.42%"
$ws.Range("E24").Value = "'-0.08%"
$ws.Range("D25").Value = "'0.001222"
$ws.Range("E25").Value = "'-0.22%"
$ws.Range("D26").Value = "'0.004152"
$ws.Range("E26").Value = "'-17.11%"
$ws.Range("D40").Value = "'0.03777"
$ws.Range("E40").Value = "'0.62%"
$ws.Range("D42").Value = "'0.003734"
$ws.Range("E42").Value = "'-34.79%"
$ws.Range("D43").Value = "'0.002489"
$ws.Range("E43").Value = "'2.56%"
$ws.Range("D44").Value = "'0.009150"
$ws.Range("E44").Value = "'7.61%"
$ws.Range("D45").Value = "'0.00005266"
$ws.Range("E45").Value = "'0.30%"
$ws.Range("E46").Value = "'-0.01%"
$ws.Range("D47").Value = "'0.1049"
$ws.Range("E47").Value = "'62.27%"
$ws.Range("E48").Value = "'-10.42%"
$ws.Range("E49").Value = "'-0.01%"
$ws.Range("E50").Value = "'-0.01%"
